# Commit: "Added stuff to sidebar"
#
# Rows 35-39 of the "2013-2022" sheet stored the last five quarters as text
# labels ("Apr-Jun 2021" ... "Apr-Jun 2022") formatted with a bold, centered
# style. Convert them to real dates (first day of each quarter) formatted
# with the same quarter-style number formats already used by the rest of
# column A (the "Jan-Mar"/"Apr-Jun"/"Jul-Sep"/"Oct-Dec" yyyy custom formats),
# so the whole column is now consistently typed as dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing quarter-date formatting (font/number format/alignment)
# down from the matching quarter earlier in the column so the new date
# cells render exactly like their peers ("Apr-Jun" yyyy, "Jul-Sep" yyyy, …)
$ws.Range("A31").Copy() | Out-Null
$ws.Range("A35").PasteSpecial(-4122) | Out-Null

$ws.Range("A32").Copy() | Out-Null
$ws.Range("A36").PasteSpecial(-4122) | Out-Null

$ws.Range("A33").Copy() | Out-Null
$ws.Range("A37").PasteSpecial(-4122) | Out-Null

$ws.Range("A30").Copy() | Out-Null
$ws.Range("A38").PasteSpecial(-4122) | Out-Null

$ws.Range("A31").Copy() | Out-Null
$ws.Range("A39").PasteSpecial(-4122) | Out-Null

# Replace the quarter-label text with the real date serials (first day of
# each quarter, Excel's 1900 date system) now that the cells carry a date
# number format: 2021-04-01, 2021-07-01, 2021-10-01, 2022-01-01, 2022-04-01.
$ws.Range("A35").Value = 44287
$ws.Range("A36").Value = 44378
$ws.Range("A37").Value = 44470
$ws.Range("A38").Value = 44562
$ws.Range("A39").Value = 44652

# Update the active selection / scroll position left by the editor.
$ws.Range("A11").Select() | Out-Null
